$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force text to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.153.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.176.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.175.60"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.725.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.173.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.173.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000117"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.818"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.610.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "316.59"
$ws.Range("D48").Style = "Normal"

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  -3.38%  "
$ws.Range("E3").Value = "  -8.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E6").Value = "  -5.02%  "
$ws.Range("E7").Value = "  -3.18%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -8.24%  "
$ws.Range("E11").Value = "  -5.23%  "
$ws.Range("E12").Value = "  -5.67%  "
$ws.Range("E13").Value = "  -8.31%  "
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("E15").Value = "  -9.46%  "
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("E17").Value = "  -5.40%  "
$ws.Range("E18").Value = "  -8.96%  "
$ws.Range("E19").Value = "  -4.46%  "
$ws.Range("E20").Value = "  -6.47%  "
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("E22").Value = "  -6.41%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  -6.51%  "
$ws.Range("E25").Value = "  -5.94%  "
$ws.Range("E26").Value = "  -6.76%  "
$ws.Range("E27").Value = "  -3.87%  "
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  -6.78%  "
$ws.Range("E32").Value = "  -5.01%  "
$ws.Range("E33").Value = "  -7.38%  "
$ws.Range("E34").Value = "  -6.31%  "
$ws.Range("E35").Value = "  -5.59%  "
$ws.Range("E36").Value = "  -8.40%  "
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("E38").Value = "  -7.74%  "
$ws.Range("E39").Value = "  -6.73%  "
$ws.Range("E40").Value = "  -6.63%  "
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("E42").Value = "  -7.00%  "
$ws.Range("E43").Value = "  -7.69%  "
$ws.Range("E44").Value = "  -7.08%  "
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("E46").Value = "  -6.95%  "
$ws.Range("E47").Value = "  -6.38%  "
$ws.Range("E48").Value = "  -8.03%  "
$ws.Range("E49").Value = "  -7.31%  "
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("E51").Value = "  -0.03%  "

Write-Host "Applied cryptos list update"